# Auto-update draw results: append the 2025-10-31 Pick 3 draw as a new
# row at the bottom of the "Results" sheet (row 45), matching the layout
# of every existing row (Date, Game, Phase, Result, InsertedAt).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.Cells.Item($ws.UsedRange.Rows.Count, 1).Row + 1

$date        = "2025-10-31"
$game        = "Pick 3"
$phase       = "251031"
$result      = "2-4-1"
$insertedAt  = "2025-10-31T21:39:16.139+04:00"

$rowRange = $ws.Range("A" + $newRow + ":E" + $newRow)

# "2025-10-31" and "251031" look like a date / a plain number, so Excel
# would silently reinterpret them unless the cells are entered as Text.
# Leading the literal with an apostrophe forces a text entry (the classic
# "quote prefix" trick), exactly like typing '2025-10-31 into the Excel
# UI; reapplying the Normal style afterwards clears the quote-prefix
# indicator so the cells end up plain text, same as the rest of the
# column.
$ws.Cells.Item($newRow, 1).Value = "'" + $date
$ws.Cells.Item($newRow, 2).Value = $game
$ws.Cells.Item($newRow, 3).Value = "'" + $phase
$ws.Cells.Item($newRow, 4).Value = $result
$ws.Cells.Item($newRow, 5).Value = $insertedAt

$rowRange.Style = "Normal"
